# Error Calculations and Plots
# - Remove the "RM 232" data row and the "SC 92" data row entirely.
# - Fill in / clear several previously-missing (or now-missing) values in
#   column C ("B" header) and column D ("C" header) as part of the error
#   (missing-data) re-imputation pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that no longer belong in the cleaned table ---
# Row 26 is "RM 232"; after it is removed, what used to be row 28 ("SC 92")
# shifts up to row 27, so we delete that row next.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# --- Column C ("B") value updates (post row-deletion row numbers) ---
$ws.Cells.Item(26, 3).Value = 10.8     # SC 5
$ws.Cells.Item(27, 3).Value = ""       # SC 101 (now missing)
$ws.Cells.Item(30, 3).Value = 11.4     # SC 120
$ws.Cells.Item(32, 3).Value = ""       # SC 193 (now missing)

# --- Column D ("C") value updates (post row-deletion row numbers) ---
$ws.Cells.Item(2, 4).Value = -13.5     # RM 2
$ws.Cells.Item(6, 4).Value = ""        # RM 21 (now missing)
$ws.Cells.Item(12, 4).Value = -14.1    # RM 81
$ws.Cells.Item(14, 4).Value = ""       # RM 90 (now missing)
$ws.Cells.Item(20, 4).Value = -14      # RM 134
$ws.Cells.Item(21, 4).Value = -14.3    # RM 135
$ws.Cells.Item(23, 4).Value = ""       # RM 140 (now missing)
$ws.Cells.Item(24, 4).Value = ""       # RM 142a (now missing)
$ws.Cells.Item(31, 4).Value = -13.7    # SC 132
$ws.Cells.Item(33, 4).Value = -14.1    # SC 232
